$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-set NumberFormat to Text for D-column cells whose new values would
# otherwise be auto-coerced to numbers by Excel (single-dot numeric strings),
# so they round-trip as text exactly like the rest of the Price column.
$forcedTextCells = @("D4","D5","D6","D7","D8","D9","D11","D14","D15","D17","D18","D19","D21","D23","D25","D26","D27","D28","D29","D30","D32","D33","D34","D35","D36","D37","D38","D39","D40","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $forcedTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated coin data cell by cell, in row order.
# Row 2
$ws.Range("D2").Value = '22.465.14'
$ws.Range("E2").Value = '  -0.08%  '

# Row 3
$ws.Range("D3").Value = '1.572.25'
$ws.Range("E3").Value = '  -0.06%  '

# Row 4
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.08%  '

# Row 5
$ws.Range("D5").Value = '1.002'

# Row 6
$ws.Range("D6").Value = '286.68'
$ws.Range("E6").Value = '  -1.85%  '

# Row 7
$ws.Range("D7").Value = '0.3658'
$ws.Range("E7").Value = '  -1.74%  '

# Row 8
$ws.Range("D8").Value = '48.12'
$ws.Range("E8").Value = '  -3.41%  '

# Row 9
$ws.Range("D9").Value = '0.3331'
$ws.Range("E9").Value = '  -2.12%  '

# Row 10
$ws.Range("E10").Value = '  -1.86%  '

# Row 11
$ws.Range("D11").Value = '0.07433'
$ws.Range("E11").Value = '  -1.56%  '

# Row 13
$ws.Range("E13").Value = '  -1.89%  '

# Row 14
$ws.Range("D14").Value = '5.985'
$ws.Range("E14").Value = '  -0.97%  '

# Row 15
$ws.Range("D15").Value = '6.914'
$ws.Range("E15").Value = '  -0.77%  '

# Row 16
$ws.Range("D16").Value = '1.571.64'
$ws.Range("E16").Value = '  +0.07%  '

# Row 17
$ws.Range("D17").Value = '0.00001108'
$ws.Range("E17").Value = '  -1.63%  '

# Row 18
$ws.Range("D18").Value = '88.01'
$ws.Range("E18").Value = '  -3.13%  '

# Row 19
$ws.Range("D19").Value = '0.06741'
$ws.Range("E19").Value = '  -0.41%  '

# Row 21
$ws.Range("D21").Value = '6.391'
$ws.Range("E21").Value = '  +1.34%  '

# Row 22
$ws.Range("E22").Value = '  +0.26%  '

# Row 23
$ws.Range("D23").Value = '12.06'
$ws.Range("E23").Value = '  -0.86%  '

# Row 24
$ws.Range("D24").Value = '22.474.78'
$ws.Range("E24").Value = '  +0.01%  '

# Row 25
$ws.Range("D25").Value = '2.383'
$ws.Range("E25").Value = '  +0.53%  '

# Row 26
$ws.Range("D26").Value = '2.608'
$ws.Range("E26").Value = '  -0.85%  '

# Row 27
$ws.Range("D27").Value = '152.15'
$ws.Range("E27").Value = '  +1.79%  '

# Row 28
$ws.Range("D28").Value = '19.57'
$ws.Range("E28").Value = '  -2.26%  '

# Row 29
$ws.Range("D29").Value = '5.017'
$ws.Range("E29").Value = '  -0.69%  '

# Row 30
$ws.Range("D30").Value = '124.20'
$ws.Range("E30").Value = '  -1.06%  '

# Row 31
$ws.Range("D31").Value = '1.747.96'
$ws.Range("E31").Value = '  +0.14%  '

# Row 32
$ws.Range("D32").Value = '1.036'
$ws.Range("E32").Value = '  -4.16%  '

# Row 33
$ws.Range("D33").Value = '6.154'
$ws.Range("E33").Value = '  -0.91%  '

# Row 34
$ws.Range("D34").Value = '1.992'
$ws.Range("E34").Value = '  -1.17%  '

# Row 35
$ws.Range("D35").Value = '9.726'
$ws.Range("E35").Value = '  -0.90%  '

# Row 36
$ws.Range("D36").Value = '0.08262'
$ws.Range("E36").Value = '  -1.42%  '

# Row 37
$ws.Range("D37").Value = '0.02428'
$ws.Range("E37").Value = '  -2.02%  '

# Row 38
$ws.Range("B38").Value = 'Algorand'
$ws.Range("C38").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D38").Value = '0.2256'
$ws.Range("E38").Value = '  -2.03%  '

# Row 39
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").Value = '0.06477'
$ws.Range("E39").Value = '  -1.00%  '

# Row 40
$ws.Range("D40").Value = '5.409'
$ws.Range("E40").Value = '  -0.70%  '

# Row 41
$ws.Range("E41").Value = '  -3.35%  '

# Row 42
$ws.Range("D42").Value = '11.29'
$ws.Range("E42").Value = '  -0.65%  '

# Row 43
$ws.Range("D43").Value = '0.6261'
$ws.Range("E43").Value = '  +0.29%  '

# Row 44
$ws.Range("D44").Value = '13.90'
$ws.Range("E44").Value = '  -1.37%  '

# Row 45
$ws.Range("D45").Value = '0.6046'
$ws.Range("E45").Value = '  +3.38%  '

# Row 46
$ws.Range("D46").Value = '3.749'
$ws.Range("E46").Value = '  -1.74%  '

# Row 47
$ws.Range("D47").Value = '2.048'

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '124.32'
$ws.Range("E48").Value = '  -4.72%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").Value = '1.221'
$ws.Range("E49").Value = '  +0.63%  '

# Row 50
$ws.Range("D50").Value = '0.07215'
$ws.Range("E50").Value = '  -1.66%  '

# Row 51
$ws.Range("D51").Value = '76.44'
$ws.Range("E51").Value = '  -0.37%  '

# Restore default (Normal) style on the forced-text cells so only the
# value/type changed, matching the original (unstyled) General-format cells.
foreach ($addr in $forcedTextCells) {
    $ws.Range($addr).Style = "Normal"
}
